# Regenerate the orders sheet with updated distance/size labels.
# Mapping applied to every text cell in the used range:
#   D80 -> D86   D51 -> D55   D64 -> D69   S30 -> S31
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.UsedRange
$arr = $rng.Value2

$rowCount = $arr.GetUpperBound(0)
$colCount = $arr.GetUpperBound(1)

# Compute replacement values up front (keeps non-string cells untouched).
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $val = $arr[$r, $c]
        if ($val -is [string]) {
            $newVal = $val -replace "D80", "D86"
            $newVal = $newVal -replace "D51", "D55"
            $newVal = $newVal -replace "D64", "D69"
            $newVal = $newVal -replace "S30", "S31"
            $arr[$r, $c] = $newVal
        }
    }
}

# Write the updated values back column by column, mirroring the
# column-major order the source data was originally generated in.
for ($c = 1; $c -le $colCount; $c++) {
    for ($r = 1; $r -le $rowCount; $r++) {
        $ws.Cells.Item($r, $c).Value2 = $arr[$r, $c]
    }
}
